$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H ("Save"), matching the style used by the other header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data values for column H (the "Save" column)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
